# Legislator property disclosure workbook — "存款" (bank deposits) sheet
# rebuild: turn the old partial/duplicated rows into a proper table with
# bank / deposit_type / currency / owner / total / property_category /
# category / date / legislator_name / legislator_id / source_file / index
# columns (commit message: "#5: cash & deposit done").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# Wipe the previous A1:G13 content but keep the header/body styles (s="1"
# on row 1 + column A, s="2" on the rest) that are already attached to
# those cells.
$ws.Range("A1:G13").ClearContents()

# ---- header row -------------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1
$ws.Range("I1:M1").Font.Bold = $true
$ws.Range("I1:M1").HorizontalAlignment = -4108
$ws.Range("I1:M1").VerticalAlignment = -4160
$ws.Range("I1:M1").Borders.LineStyle = 1

# ---- data rows ----------------------------------------------------------
# index -> (A, bank, deposit_type, currency, owner, total, property_category,
#           category, date, legislator_name, legislator_id, source_file, index)
# the date column is text like "2012-04-27" in the original file, not a
# real date — a leading apostrophe keeps Excel's autodetect from turning it
# into a date serial number.
$d = "'2012-04-27"

$rows = @(
  @(45, "第一商業銀行",     "活期存款", "新臺幣", "李貴敏", 1856851,   "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 45),
  @(46, "第商業銀行",       "活期存款", "新臺幣", "李貴敏", 83562,     "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 46),
  @(47, "第一商業銀行",     "活期存款", "新臺幣", "李貴敏", 101986,    "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 47),
  @(48, "第一商業銀行",     "活期存款", "新臺幣", "李貴敏", 8222,      "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 48),
  @(49, "中國信託商業銀行", "活期存款", "新臺幣", "李貴敏", 86355,     "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 49),
  @(50, "安泰商業銀行",     "活期存款", "新臺幣", "李貴敏", 44956,     "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 50),
  @(51, "兆豐國際商業銀行", "活期存款", "新臺幣", "李貴敏", 1000,      "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 51),
  @(52, "台北富邦商業銀行", "活期存款", "新臺幣", "李貴敏", 207092,    "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 52),
  @(53, "臺灣銀行",         "活期存款", "新臺幣", "李貴敏", 930984,    "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 53),
  @(54, "中華郵政股份有限公司", "活期存款", "新臺幣", "李貴敏", 852053, "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 54),
  @(55, "第一商業銀行",     "活期存款", "美金",   "李貴敏", 1.18,      "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 55),
  @(56, "台北富邦商業銀行", "活期存款", "美金",   "李貴敏", 457974.81, "deposit", "normal", $d, "李貴敏", 1739, "tmp59331", 56)
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value  = $row[0]   # A  index (display)
  $ws.Cells.Item($r, 2).Value  = $row[1]   # B  bank
  $ws.Cells.Item($r, 3).Value  = $row[2]   # C  deposit_type
  $ws.Cells.Item($r, 4).Value  = $row[3]   # D  currency
  $ws.Cells.Item($r, 5).Value  = $row[4]   # E  owner
  $ws.Cells.Item($r, 6).Value  = $row[5]   # F  total
  $ws.Cells.Item($r, 7).Value  = $row[6]   # G  property_category
  $ws.Cells.Item($r, 8).Value  = $row[7]   # H  category
  $ws.Cells.Item($r, 9).Value  = $row[8]   # I  date
  $ws.Cells.Item($r, 10).Value = $row[9]   # J  legislator_name
  $ws.Cells.Item($r, 11).Value = $row[10]  # K  legislator_id
  $ws.Cells.Item($r, 12).Value = $row[11]  # L  source_file
  $ws.Cells.Item($r, 13).Value = $row[12]  # M  index

  $r = $r + 1
}
